$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = 'Datos actualizados a 27 de Abril de 2020 a las 20:52'
$ws.Range("B4").Value = 59421
$ws.Range("C4").Value = 35565
$ws.Range("D4").Value = 15870
$ws.Range("E4").Value = 7986
$ws.Range("B5").Value = 47755
$ws.Range("C5").Value = 17897
$ws.Range("D5").Value = 25159
$ws.Range("E5").Value = 4699
$ws.Range("A6").Value = 'Castilla y Leon'
$ws.Range("B6").Value = 16404
$ws.Range("C6").Value = 6272
$ws.Range("D6").Value = 8442
$ws.Range("E6").Value = 1690
$ws.Range("A7").Value = 'Castilla-La Mancha'
$ws.Range("B7").Value = 15664
$ws.Range("C7").Value = 5282
$ws.Range("D7").Value = 8017
$ws.Range("E7").Value = 2365
$ws.Range("B8").Value = 12513
$ws.Range("C8").Value = 9840
$ws.Range("D8").Value = 1432
$ws.Range("E8").Value = 1241
$ws.Range("B9").Value = 11852
$ws.Range("C9").Value = 4890
$ws.Range("D9").Value = 5805
$ws.Range("E9").Value = 1157
$ws.Range("B10").Value = 9238
$ws.Range("C10").Value = 1816
$ws.Range("D10").Value = 7017
$ws.Range("E10").Value = 405
$ws.Range("A13").Value = 'Valencia/Valencia'
$ws.Range("B13").Value = 5131
$ws.Range("C13").Value = 2194
$ws.Range("D13").Value = 2583
$ws.Range("E13").Value = 515
$ws.Range("A14").Value = 'Aragon'
$ws.Range("B14").Value = 4985
$ws.Range("C14").Value = 2010
$ws.Range("D14").Value = 2249
$ws.Range("E14").Value = 726
$ws.Range("A15").Value = 'Navarra'
$ws.Range("B15").Value = 4733
$ws.Range("C15").Value = 1918
$ws.Range("D15").Value = 2383
$ws.Range("E15").Value = 432
$ws.Range("A16").Value = 'Toledo'
$ws.Range("B16").Value = 3938
$ws.Range("C16").Value = 4178
$ws.Range("D16").Value = 10597
$ws.Range("E16").Value = 504
$ws.Range("A17").Value = 'La Rioja'
$ws.Range("B17").Value = 3892
$ws.Range("C17").Value = 2049
$ws.Range("D17").Value = 1523
$ws.Range("E17").Value = 320
$ws.Range("A21").Value = 'Araba/Alava'
$ws.Range("B21").Value = 3241
$ws.Range("C21").Value = 7124
$ws.Range("D21").Value = 4423
$ws.Range("E21").Value = 318
$ws.Range("A22").Value = 'Valladolid'
$ws.Range("B22").Value = 3154
$ws.Range("C22").Value = 1070
$ws.Range("D22").Value = 1824
$ws.Range("E22").Value = 260
$ws.Range("A23").Value = 'Extremadura'
$ws.Range("B23").Value = 2749
$ws.Range("C23").Value = 1652
$ws.Range("D23").Value = 669
$ws.Range("E23").Value = 428
$ws.Range("A25").Value = 'Malaga'
$ws.Range("B25").Value = 2531
$ws.Range("C25").Value = 869
$ws.Range("D25").Value = 1439
$ws.Range("E25").Value = 223
$ws.Range("A26").Value = 'Segovia'
$ws.Range("B26").Value = 2406
$ws.Range("C26").Value = 656
$ws.Range("D26").Value = 1578
$ws.Range("E26").Value = 172
$ws.Range("A27").Value = 'Leon'
$ws.Range("B27").Value = 2403
$ws.Range("C27").Value = 1076
$ws.Range("D27").Value = 1024
$ws.Range("E27").Value = 303
$ws.Range("A28").Value = 'Gipuzkoa/Guipuzcoa'
$ws.Range("B28").Value = 2342
$ws.Range("C28").Value = 7124
$ws.Range("D28").Value = 4423
$ws.Range("E28").Value = 212
$ws.Range("A29").Value = 'Sevilla'
$ws.Range("B29").Value = 2329
$ws.Range("C29").Value = 459
$ws.Range("D29").Value = 1658
$ws.Range("A30").Value = 'Asturias'
$ws.Range("B30").Value = 2254
$ws.Range("C30").Value = 759
$ws.Range("D30").Value = 1242
$ws.Range("E30").Value = 253
$ws.Range("A31").Value = 'Caceres'
$ws.Range("B31").Value = 2220
$ws.Range("C31").Value = 422
$ws.Range("D31").Value = 1482
$ws.Range("E31").Value = 316
$ws.Range("A32").Value = 'Gran Canaria'
$ws.Range("B32").Value = 2178
$ws.Range("C32").Value = 1047
$ws.Range("D32").Value = 1000
$ws.Range("E32").Value = 131
$ws.Range("A33").Value = 'Cantabria'
$ws.Range("B33").Value = 2104
$ws.Range("C33").Value = 1241
$ws.Range("D33").Value = 677
$ws.Range("E33").Value = 186
$ws.Range("A36").Value = 'Burgos'
$ws.Range("B36").Value = 1567
$ws.Range("C36").Value = 642
$ws.Range("D36").Value = 757
$ws.Range("E36").Value = 168
$ws.Range("A37").Value = 'Pontevedra'
$ws.Range("B37").Value = 1536
$ws.Range("C37").Value = 333
$ws.Range("D37").Value = 1411
$ws.Range("E37").Value = 30
$ws.Range("A38").Value = 'Murcia'
$ws.Range("B38").Value = 1474
$ws.Range("C38").Value = 990
$ws.Range("D38").Value = 356
$ws.Range("E38").Value = 128
$ws.Range("B59").Value = 110
$ws.Range("C59").Value = 87
$ws.Range("D59").Value = 21
